# Apply updated experiment values (non-convex experiments, all but the 5th).
# All target cells in this workbook store their contents as plain text
# (shared strings), even when the text looks like a number. Assigning a
# numeric-looking string straight to .Value would make Excel auto-convert
# it into a real number cell, so we briefly force Text number format,
# assign the value, then clear the formatting again to avoid leaving any
# visible style change behind.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$cellAddr, [string]$val)
    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# NOTE: worksheet names are looked up case-insensitively, and this workbook
# has two sheets whose names differ only by case ("Vector_bf" / "Vector_BF").
# Using Worksheets.Item(<name>) would therefore resolve both to the same
# sheet, so every sheet below is addressed by its (unambiguous) 1-based
# position instead, per the order in the Sheets collection:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# --- Restricciones_del_lider ---
$wsLider = $wb.Worksheets.Item(2)
Set-TextValue $wsLider "A2" "-2.1 + x"
Set-TextValue $wsLider "B2" "1.1"
Set-TextValue $wsLider "D2" "0.74"
Set-TextValue $wsLider "A3" "2.1 - x"
Set-TextValue $wsLider "B3" "-3.1"
Set-TextValue $wsLider "D3" "0.27"

# --- Restricciones_del_follower ---
$wsFollower = $wb.Worksheets.Item(3)
Set-TextValue $wsFollower "A2" "-3.3000000000000003 + y"
Set-TextValue $wsFollower "B2" "2.3000000000000003"
Set-TextValue $wsFollower "D2" "0.22"
Set-TextValue $wsFollower "E2" "-3.5"
Set-TextValue $wsFollower "F2" "-6.0"
Set-TextValue $wsFollower "A3" "3.3000000000000007 - y"
Set-TextValue $wsFollower "B3" "-4.300000000000001"
Set-TextValue $wsFollower "D3" "0.66"
Set-TextValue $wsFollower "E3" "-4.9"
Set-TextValue $wsFollower "F3" "-6.800000000000001"

# --- Punto_modificado ---
$wsPunto = $wb.Worksheets.Item(4)
Set-TextValue $wsPunto "A2" "2.1"
Set-TextValue $wsPunto "B2" "3.3000000000000003"

# --- Vector_bf ---
$wsBf = $wb.Worksheets.Item(5)
Set-TextValue $wsBf "A2" "-8.350000000000003"

# --- Vector_BF ---
$wsBF = $wb.Worksheets.Item(6)
Set-TextValue $wsBF "A2" "-4.17"
Set-TextValue $wsBF "A3" "-8.0"
